$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test data values
$ws.Range("A2").Value = "orderVegetables"
$ws.Range("A3").Value = "orderFruits"
$ws.Range("B2").Value = "Onion,Cucumber,Brocolli,Musk Melon"
$ws.Range("C2").Value = "3,2,2,10"

# Apply new formatting to A2 (Courier New 10pt, black, vertically centered, no border)
$ws.Range("A2").Font.Name = "Courier New"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.Color = 0
$ws.Range("A2").VerticalAlignment = -4108
$ws.Range("A2").Borders.LineStyle = -4142

# Update the active selection
$ws.Range("D7").Select()
